$wb = $excel.ActiveWorkbook

# --- Summary sheet: refresh headline metrics after trade #95 closed ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.68               # Current Capital
$wsSummary.Range("B4").Value = -0.31                 # Total P&L $
$wsSummary.Range("B5").Value = -0.07000000000000001  # Total P&L %
$wsSummary.Range("B6").Value = 95                    # Total Trades
$wsSummary.Range("B7").Value = 39                    # Winning Trades
$wsSummary.Range("B9").Value = 41.05                 # Win Rate %

# --- Strategy Status sheet: MarketMaking strategy row (row 4) ---
$wsStrategy = $wb.Worksheets.Item("Strategy Status")
$wsStrategy.Range("C4").Value = 99.68000000000001    # Capital
$wsStrategy.Range("D4").Value = 95                   # Trades
$wsStrategy.Range("E4").Value = -0.31                # P&L $
$wsStrategy.Range("F4").Value = -0.32                # P&L %
$wsStrategy.Range("G4").Value = 41.05                # Win Rate %

# --- New trade row (#95, 1-based index 96 with header) appended to both the
#     "All Trades" log and the per-strategy "MarketMaking" log ---
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 96

    $ws.Cells.Item($row, 1).Value  = 95
    # Leading apostrophe forces these to be stored as text (matching the
    # existing Date/Time columns) instead of being auto-parsed into date
    # serial numbers.
    $ws.Cells.Item($row, 2).Value  = "'2026-02-17"
    $ws.Cells.Item($row, 3).Value  = "'09:09:23"
    $ws.Cells.Item($row, 4).Value  = "MarketMaking"
    $ws.Cells.Item($row, 5).Value  = "UP"
    $ws.Cells.Item($row, 6).Value  = 0.86
    $ws.Cells.Item($row, 7).Value  = 0.95
    $ws.Cells.Item($row, 8).Value  = "CLOSED"
    $ws.Cells.Item($row, 9).Value  = 10.4651
    $ws.Cells.Item($row, 10).Value = 0.09
    $ws.Cells.Item($row, 11).Value = 99.68000000000001
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.1

    # Drop the quote-prefix style picked up from the leading apostrophes so
    # the new row's cells stay on the workbook's default (unstyled) format.
    $ws.Range("A96:Q96").Style = "Normal"
}
